$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# New header cell: "Compliant" column, appended as column E
$ws.Cells.Item(1, 5).Value = "Compliant"

# Matching new column-width entry for column F (width stored as 17 in xml;
# ColumnWidth property is offset from the stored width by the default
# 5/6-character padding, so subtract that to land exactly on 17).
$ws.Columns.Item(6).ColumnWidth = 17 - 0.8333333333333334
